{"js": "// Lekcja 11 - Debugowanie: fix the typo \"do wybranej linii\" -> \"od wybranej linii\"\n// in the phrase \"Kontynuowanie dzia\u0142ania programu do/od wybranej linii\".\n//\n// This phrase appears three times in the document:\n//   1) As the actual heading text (wrapped by the \"_Ref86091298\" /\n//      \"_Hlk86091493\" bookmarks used for cross references).\n//   2) and 3) As the cached/displayed text of the two \"REF _Ref86091298\"\n//      cross-reference fields elsewhere in the document.\n//\n// We search for the narrower substring \" do wybranej linii\" (with the\n// leading space) so that only the trailing part of the phrase is touched;\n// this leaves the preceding \"Kontynuowanie dzia\u0142ania programu\" text (and\n// therefore the nested \"_Hlk86091493\" bookmark that wraps only that part)\n// untouched.\n\nconst body = context.document.body;\nconst results = body.search(\" do wybranej linii\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\" od wybranej linii\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Lekcja 11 - Debugowanie: fix the typo \"do wybranej linii\" -> \"od wybranej linii\"\n# in the phrase \"Kontynuowanie dzia\u0142ania programu do/od wybranej linii\".\n#\n# This phrase appears three times in the document:\n#   1) As the actual heading text (wrapped by the \"_Ref86091298\" /\n#      \"_Hlk86091493\" bookmarks used for cross references).\n#   2) and 3) As the cached/displayed text of the two \"REF _Ref86091298\"\n#      cross-reference fields elsewhere in the document.\n#\n# We search/replace the narrower string \" do wybranej linii\" (with the\n# leading space) so only the trailing part of the phrase is touched; this\n# leaves the preceding \"Kontynuowanie dzia\u0142ania programu\" run (and so the\n# nested \"_Hlk86091493\" bookmark that wraps only that part) untouched, and\n# replaces every matching occurrence in the document body.\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \" do wybranej linii\"\n$find.Replacement.Text = \" od wybranej linii\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
